# regen sval data to filter save games
# Replaces the B:E (and derived G) score columns on rows 2-30 of Sheet1
# with the regenerated values. Column A (dates) and F (win flag) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (B, C, D, E) new values. G (sum) is recomputed as B+C+D+E.
$rows = @{
    2  = @(3.272327238179451,  1.626987699542094,  0.7210945179870265, 0.5333859586016987)
    3  = @(1.445647641019636,  1.626987699542094,  3.223369029078222,  0.5333859586016987)
    4  = @(3.272327238179451,  1.626987699542094,  3.223369029078222,  0.5333859586016987)
    5  = @(0.1169995834814548, 0.04103571897497393,3.223369029078222,  0.5333859586016987)
    6  = @(3.272327238179451,  1.626987699542094,  0.7210945179870265, 0.5333859586016987)
    7  = @(1.445647641019636,  1.626987699542094,  0.1496068669990043, 0.5333859586016987)
    8  = @(0.2881169905109251, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987)
    9  = @(3.272327238179451,  1.626987699542094,  3.223369029078222,  0.5333859586016987)
    10 = @(1.445647641019636,  1.626987699542094,  0.7210945179870265, 0.5333859586016987)
    11 = @(1.445647641019636,  0.3048912486333797, 0.7210945179870265, 0.5333859586016987)
    12 = @(3.272327238179451,  1.626987699542094,  0.1496068669990043, 0.5333859586016987)
    13 = @(1.445647641019636,  1.626987699542094,  3.223369029078222,  0.5333859586016987)
    14 = @(3.272327238179451,  1.626987699542094,  0.7210945179870265, 0.5333859586016987)
    15 = @(0.2881169905109251, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987)
    16 = @(0.6545652718822623, 9.983522426115931,  0.1496068669990043, 13.86384647080068)
    17 = @(1.445647641019636,  0.3048912486333797, 0.7210945179870265, 0.5333859586016987)
    18 = @(3.272327238179451,  1.626987699542094,  3.223369029078222,  0.5333859586016987)
    19 = @(1.445647641019636,  1.626987699542094,  18.71679738969934,  13.86384647080068)
    20 = @(3.272327238179451,  1.626987699542094,  18.71679738969934,  0.5333859586016987)
    21 = @(0.6545652718822623, 1.626987699542094,  0.1496068669990043, 0.5333859586016987)
    22 = @(3.272327238179451,  1.626987699542094,  0.1496068669990043, 0.5333859586016987)
    23 = @(1.445647641019636,  1.626987699542094,  0.7210945179870265, 0.5333859586016987)
    24 = @(3.272327238179451,  1.626987699542094,  3.223369029078222,  0.5333859586016987)
    25 = @(3.272327238179451,  1.626987699542094,  0.1496068669990043, 0.5333859586016987)
    26 = @(1.445647641019636,  1.626987699542094,  3.223369029078222,  0.5333859586016987)
    27 = @(3.272327238179451,  1.626987699542094,  0.7210945179870265, 0.5333859586016987)
    28 = @(0.2881169905109251, 0.3048912486333797, 3.223369029078222,  0.5333859586016987)
    29 = @(0.6545652718822623, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987)
    30 = @(0.2881169905109251, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]
    $g = $b + $c + $d + $e

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 7).Value = $g
}
